$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44418
$ws.Range("L2").Value = 'Especial'
$ws.Range("M2").Value = 100
$ws.Range("N2:P2").Value = 8000
$ws.Range("Q2").Value = '$/caja 15 kilos granel'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 533
$ws.Range("T2").Value = 15

# Row 4
$ws.Range("D4").Value = 44217
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 55
$ws.Range("N4:P4").Value = 18000
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 1000
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44966
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 4
$ws.Range("N5:P5").Value = 250000
$ws.Range("Q5").Value = '$/bins (400 kilos)'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 625
$ws.Range("T5").Value = 400

# Row 6
$ws.Range("D6").Value = 44966
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 80
$ws.Range("N6:P6").Value = 15000
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 833
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("D7").Value = 44208
$ws.Range("L7").Value = 'Especial'
$ws.Range("M7").Value = 70
$ws.Range("N7:P7").Value = 24000
$ws.Range("Q7").Value = '$/caja 15 kilos granel'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 1600
$ws.Range("T7").Value = 15

# Row 8
$ws.Range("D8").Value = 44495
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 50
$ws.Range("N8:P8").Value = 24000
$ws.Range("Q8").Value = '$/bandeja 10 kilos'
$ws.Range("R8").Value = 'China'
$ws.Range("S8").Value = 2400
$ws.Range("T8").Value = 10

# Row 9
$ws.Range("D9").Value = 44264
$ws.Range("L9").Value = 'Calibre 100'
$ws.Range("M9").Value = 50
$ws.Range("N9:P9").Value = 20000
$ws.Range("Q9").Value = '$/caja 18 kilos embalada'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 1111
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44427
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 55
$ws.Range("N10:P10").Value = 7000
$ws.Range("Q10").Value = '$/caja 15 kilos granel'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 467
$ws.Range("T10").Value = 15

# Row 13
$ws.Range("D13").Value = 44411
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 210
$ws.Range("N13:P13").Value = 8000
$ws.Range("Q13").Value = '$/bandeja 8 kilos'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 1000
$ws.Range("T13").Value = 8
